$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 185
$ws.Cells.Item(184, 1).Copy()
$ws.Cells.Item(185, 1).PasteSpecial(-4122)
$ws.Cells.Item(185, 1).Value2 = 45615
$ws.Cells.Item(185, 2).Value2 = 817.0224939162
$ws.Cells.Item(185, 3).Value2 = 215.458547444
$ws.Cells.Item(185, 9).Value2 = 425.462993853
$ws.Cells.Item(185, 11).Value2 = 96.83497567138501
$ws.Cells.Item(185, 14).Value2 = 35.13320051136
$ws.Cells.Item(185, 15).Value2 = 1.047870285
$ws.Cells.Item(185, 17).Value2 = 0.0000047232
$ws.Cells.Item(185, 21).Value2 = 269.1626541499384
$ws.Cells.Item(185, 26).Value2 = 1795.073932935398

# Row 186
$ws.Cells.Item(185, 1).Copy()
$ws.Cells.Item(186, 1).PasteSpecial(-4122)
$ws.Cells.Item(186, 1).Value2 = 45616
$ws.Cells.Item(186, 2).Value2 = 834.5095995168001
$ws.Cells.Item(186, 3).Value2 = 212.8605144095
$ws.Cells.Item(186, 9).Value2 = 421.613541188
$ws.Cells.Item(186, 11).Value2 = 94.505111595081
$ws.Cells.Item(186, 14).Value2 = 32.30438416800001
$ws.Cells.Item(186, 15).Value2 = 1.0306583133
$ws.Cells.Item(186, 17).Value2 = 0.000004502399999999999
$ws.Cells.Item(186, 21).Value2 = 253.299455901558
$ws.Cells.Item(186, 26).Value2 = 1770.852056683318

# Row 187
$ws.Cells.Item(186, 1).Copy()
$ws.Cells.Item(187, 1).PasteSpecial(-4122)
$ws.Cells.Item(187, 1).Value2 = 45617
$ws.Cells.Item(187, 2).Value2 = 870.1831993536
$ws.Cells.Item(187, 3).Value2 = 232.6796166935
$ws.Cells.Item(187, 9).Value2 = 459.015897547
$ws.Cells.Item(187, 11).Value2 = 96.349587322155
$ws.Cells.Item(187, 14).Value2 = 35.88987978016
$ws.Cells.Item(187, 15).Value2 = 1.0586979942
$ws.Cells.Item(187, 17).Value2 = 0.000004968
$ws.Cells.Item(187, 21).Value2 = 263.2779193158618
$ws.Cells.Item(187, 26).Value2 = 1733.401309555102

# Row 188
$ws.Cells.Item(187, 1).Copy()
$ws.Cells.Item(188, 1).PasteSpecial(-4122)
$ws.Cells.Item(188, 1).Value2 = 45618
$ws.Cells.Item(188, 2).Value2 = 875.27133576
$ws.Cells.Item(188, 3).Value2 = 230.736118803
$ws.Cells.Item(188, 9).Value2 = 459.96483239
$ws.Cells.Item(188, 11).Value2 = 99.40753392230401
$ws.Cells.Item(188, 14).Value2 = 36.0295744144
$ws.Cells.Item(188, 15).Value2 = 1.07766351
$ws.Cells.Item(188, 17).Value2 = 0.000004953599999999999
$ws.Cells.Item(188, 21).Value2 = 263.6617063702581
$ws.Cells.Item(188, 26).Value2 = 1785.245440840804

# Row 189
$ws.Cells.Item(188, 1).Copy()
$ws.Cells.Item(189, 1).PasteSpecial(-4122)
$ws.Cells.Item(189, 1).Value2 = 45619
$ws.Cells.Item(189, 2).Value2 = 864.4769244720001
$ws.Cells.Item(189, 3).Value2 = 235.3213316285
$ws.Cells.Item(189, 9).Value2 = 456.366041759
$ws.Cells.Item(189, 11).Value2 = 102.90233003676
$ws.Cells.Item(189, 14).Value2 = 37.74083368384
$ws.Cells.Item(189, 15).Value2 = 1.1073205374
$ws.Cells.Item(189, 17).Value2 = 0.000004917599999999999
$ws.Cells.Item(189, 21).Value2 = 331.4640859802711
$ws.Cells.Item(189, 26).Value2 = 1706.431181959036

Write-Host "Applied rows 185-189"
